# Edit script: update keyword lists across weekday sheets, add search-result
# columns on Thursday, and refresh view/selection state (iron_native COM).
$wb = $excel.ActiveWorkbook

# ---- Monday ----
$ws = $wb.Worksheets.Item("Monday")
$ws.Range("A6").Value = "Machine "
$ws.Range("A7").Value = "Climate "
$ws.Range("A8").Value = "Cute cat "
$ws.Range("A9").Value = "Best "
$ws.Range("A10").Value = "Web "
$ws.Range("A11").Value = "Top travel "
$ws.Range("A12").Value = "Funny "
$ws.Range("A13").Value = "Cloud "
$ws.Range("A14").Value = "IELTS "
$ws.Columns.Item(1).ColumnWidth = 28.83333333333333
$ws.Columns.Item(2).ColumnWidth = 44.66666666666666
$ws.Columns.Item(3).ColumnWidth = 42.0
$ws.Range("A1:C14").Select()

# ---- Tuesday ----
$ws = $wb.Worksheets.Item("Tuesday")
$ws.Range("A2").Value = "Dhaka"
$ws.Range("A3").Value = "University"
$ws.Range("A4").Value = "Cricket"
$ws.Range("A5").Value = "Bombay"
$ws.Range("A6").Value = "Machine "
$ws.Range("A7").Value = "Climate "
$ws.Range("A8").Value = "Cute cat "
$ws.Range("A9").Value = "Best "
$ws.Range("A10").Value = "Web "
$ws.Range("A11").Value = "Top travel "
$ws.Range("A12").Value = "Funny "
$ws.Range("A13").Value = "Cloud "
$ws.Range("A14").Value = "IELTS "
$ws.Range("A1:C14").Select()

# ---- Wednesday ----
$ws = $wb.Worksheets.Item("Wednesday")
$ws.Range("A2").Value = "Dhaka"
$ws.Range("A3").Value = "University"
$ws.Range("A4").Value = "Cricket"
$ws.Range("A5").Value = "Bombay"
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("A6").Value = "Machine "
$ws.Range("A7").Value = "Climate "
$ws.Range("A8").Value = "Cute cat "
$ws.Range("A9").Value = "Best "
$ws.Range("A10").Value = "Web "
$ws.Range("A11").Value = "Top travel "
$ws.Range("A12").Value = "Funny "
$ws.Range("A13").Value = "Cloud "
$ws.Range("A14").Value = "IELTS "
$ws.Range("A1:C14").Select()

# ---- Thursday ----
$ws = $wb.Worksheets.Item("Thursday")
$ws.Range("A2").Value = "Dhaka"
$ws.Range("A3").Value = "University"
$ws.Range("A4").Value = "Cricket"
$ws.Range("A5").Value = "Bombay"
$ws.Range("A6").Value = "Machine "
$ws.Range("A7").Value = "Climate "
$ws.Range("A8").Value = "Cute cat "
$ws.Range("A9").Value = "Best "
$ws.Range("A10").Value = "Web "
$ws.Range("A11").Value = "Top travel "
$ws.Range("A12").Value = "Funny "
$ws.Range("A13").Value = "Cloud "
$ws.Range("A14").Value = "IELTS "
$ws.Range("B2").Value = "বিশ্ববিদ্যালয়, ঢাকা, বাংলাদেশ"
$ws.Range("C2").Value = "ঢাকা"
$ws.Range("B3").Value = "দ্য ইন্টারন্যাশনাল ইউনিভার্সিটি অফ স্কলারস — বেসরকারি বিশ্ববিদ্যালয়, ঢাকা, বাংলাদেশ"
$ws.Range("C3").Value = "university"
$ws.Range("B4").Value = "আইসিসি ক্রিকেট বিশ্ব কাপ — ক্রিকেট লীগ"
$ws.Range("C4").Value = "crickex"
$ws.Range("B5").Value = "বোম্বে সুইটস এন্ড চানাচুর · chawlk, Circuler Rd, ঢাকা"
$ws.Range("C5").Value = "Bombay"
$ws.Range("B6").Value = "মাচাইন গুন কেলি — আমেরিকান র‍্যাপার"
$ws.Range("C6").Value = "machine gun"
$ws.Range("B7").Value = "climate change paragraph 150 words pdf"
$ws.Range("C7").Value = "climate change"
$ws.Range("B8").Value = "cute cat profile picture"
$ws.Range("C8").Value = "cute cat pic"
$ws.Range("B9").Value = "best football player in the world"
$ws.Range("C9").Value = "Best Buy"
$ws.Range("B10").Value = "web push notifications"
$ws.Range("C10").Value = "web do"
$ws.Range("B11").Value = "top travel agencies in the world"
$ws.Range("C11").Value = "top travel movies"
$ws.Range("B12").Value = "funny birthday wishes for best friend"
$ws.Range("C12").Value = "funny pic"
$ws.Range("B13").Value = "cloud meaning in bengali"
$ws.Range("C13").Value = "cloud ai"
$ws.Range("B14").Value = "ielts listening practice"
$ws.Range("C14").Value = "ielts liz"
$ws.Range("A1:C14").Select()

# ---- Friday ----
$ws = $wb.Worksheets.Item("Friday")
$ws.Range("A2").Value = "Dhaka"
$ws.Range("A3").Value = "University"
$ws.Range("A4").Value = "Cricket"
$ws.Range("A5").Value = "Bombay"
$ws.Range("A6").Value = "Machine "
$ws.Range("A7").Value = "Climate "
$ws.Range("A8").Value = "Cute cat "
$ws.Range("A9").Value = "Best "
$ws.Range("A10").Value = "Web "
$ws.Range("A11").Value = "Top travel "
$ws.Range("A12").Value = "Funny "
$ws.Range("A13").Value = "Cloud "
$ws.Range("A14").Value = "IELTS "
$ws.Range("A1:C14").Select()

# ---- Saturday ----
$ws = $wb.Worksheets.Item("Saturday")
$ws.Range("A2").Value = "Dhaka"
$ws.Range("A3").Value = "University"
$ws.Range("A4").Value = "Cricket"
$ws.Range("A5").Value = "Bombay"
$ws.Range("A6").Value = "Machine "
$ws.Range("A7").Value = "Climate "
$ws.Range("A8").Value = "Cute cat "
$ws.Range("A9").Value = "Best "
$ws.Range("A10").Value = "Web "
$ws.Range("A11").Value = "Top travel "
$ws.Range("A12").Value = "Funny "
$ws.Range("A13").Value = "Cloud "
$ws.Range("A14").Value = "IELTS "
$ws.Range("A1:C14").Select()

# ---- Sunday ----
$ws = $wb.Worksheets.Item("Sunday")
$ws.Range("A2").Value = "Dhaka"
$ws.Range("A3").Value = "University"
$ws.Range("A4").Value = "Cricket"
$ws.Range("A5").Value = "Bombay"
$ws.Range("A6").Value = "Machine "
$ws.Range("A7").Value = "Climate "
$ws.Range("A8").Value = "Cute cat "
$ws.Range("A9").Value = "Best "
$ws.Range("A10").Value = "Web "
$ws.Range("A11").Value = "Top travel "
$ws.Range("A12").Value = "Funny "
$ws.Range("A13").Value = "Cloud "
$ws.Range("A14").Value = "IELTS "
$ws.Range("A1:C14").Select()

# Sunday is the active tab in the saved workbook view
$wb.Worksheets.Item("Sunday").Activate()

